# Update "想去人数" (F) / "最低票价" (G) figures to the latest scrape values.
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types, a superset that
# duplicates the same rows) both carry the same underlying event data, so
# every change is mirrored across both sheets. Sheet "演出" (shows) only has
# two of its rows touched, and those two rows are also mirrored inside
# "全部类型".

$wb = $excel.ActiveWorkbook

# ---- 展览 (exhibitions) ----
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 5941
$wsExpo.Range("G2").Value = "已售罄"

$wsExpo.Range("F3").Value = 563
$wsExpo.Range("F4").Value = 1151
$wsExpo.Range("F5").Value = 1082
$wsExpo.Range("F6").Value = 860
$wsExpo.Range("F7").Value = 93
$wsExpo.Range("F11").Value = 38
$wsExpo.Range("F13").Value = 2120
$wsExpo.Range("F14").Value = 1536
$wsExpo.Range("F15").Value = 1186
$wsExpo.Range("F18").Value = 463
$wsExpo.Range("F19").Value = 689
$wsExpo.Range("F20").Value = 249
$wsExpo.Range("F23").Value = 526
$wsExpo.Range("F24").Value = 3911
$wsExpo.Range("F27").Value = 114
$wsExpo.Range("F28").Value = 176
$wsExpo.Range("F30").Value = 566
$wsExpo.Range("F36").Value = 878
$wsExpo.Range("F38").Value = 79
$wsExpo.Range("F40").Value = 97

# ---- 演出 (shows) ----
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("F3").Value = 764
$wsShow.Range("F5").Value = 416

# ---- 全部类型 (all types) ----
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 5941
$wsAll.Range("G2").Value = "已售罄"

$wsAll.Range("F3").Value = 563
$wsAll.Range("F4").Value = 1151
$wsAll.Range("F6").Value = 764
$wsAll.Range("F7").Value = 1082
$wsAll.Range("F8").Value = 860
$wsAll.Range("F10").Value = 416
$wsAll.Range("F11").Value = 93
$wsAll.Range("F15").Value = 38
$wsAll.Range("F18").Value = 2120
$wsAll.Range("F19").Value = 1536
$wsAll.Range("F20").Value = 1186
$wsAll.Range("F23").Value = 463
$wsAll.Range("F25").Value = 689
$wsAll.Range("F26").Value = 249
$wsAll.Range("F29").Value = 526
$wsAll.Range("F30").Value = 3911
$wsAll.Range("F33").Value = 114
$wsAll.Range("F34").Value = 176
$wsAll.Range("F36").Value = 566
$wsAll.Range("F42").Value = 878
$wsAll.Range("F44").Value = 79
$wsAll.Range("F46").Value = 97
